# Complete tests for Accounts module
#
# - Mark "Create new account with obligatory fields" (row 7) as Done.
# - Rename "Edit existing" (row 9) to "Edit existing account" and mark it Done.
# - Insert a new Accounts test case row (row 11): "Check error messages on
#   empty acc creation", marked Done, type Regres. This pushes every
#   following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status column's green "Done" fill, taken from the existing highlighted
# cells (RGB 146, 208, 80 == 0x92D050, encoded BGR-packed for .Color).
$doneColor = 5296274

# Insert a fresh row right after row 10 ("Delete new account"); everything
# from the old row 11 ("Contacts") onward shifts down to make room for the
# new Accounts test case.
$ws.Rows("11").Insert()

# Row 7 now has its Status marked Done.
$ws.Range("D7").Value = "Done"
$ws.Range("D7").Interior.Color = $doneColor

# Row 9: rename the test case and mark it Done.
$ws.Range("B9").Value = "Edit existing account"
$ws.Range("D9").Value = "Done"
$ws.Range("D9").Interior.Color = $doneColor

# New row 11: the additional Accounts test case.
$ws.Range("B11").Value = "Check error messages on empty acc creation"
$ws.Range("D11").Value = "Done"
$ws.Range("D11").Interior.Color = $doneColor
$ws.Range("E11").Value = "Regres"

# Matches the saved selection left behind in the source file.
$ws.Range("D10").Select()
